{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Change 1 -----------------------------------------------------------\n// The paragraph \"2.06.2022: Implementierung von weitern Funktionen der\n// GameEngine:\" is currently split across two <w:r> runs (\"2.06.2022\" and\n// \": Implementierung ...\"). Re-insert the full text as a single\n// InsertLocation.replace call so it collapses back down to one run.\nconst mergedText =\n  \"2.06.2022: Implementierung von weitern Funktionen der GameEngine:\";\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === mergedText) {\n    items[i].getRange().insertText(mergedText, Word.InsertLocation.replace);\n    break;\n  }\n}\nawait context.sync();\n\n// --- Change 2 -----------------------------------------------------------\n// Of the two empty trailing paragraphs, the first one (immediately after\n// \"Beginn des Texture-Designs\") receives the new text \"6.06.2022: \"; the\n// very last (empty) paragraph of the document is left untouched.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items2 = paragraphs.items;\nfor (let i = 1; i < items2.length - 1; i++) {\n  const prevText = items2[i - 1].text.trim();\n  const curText = items2[i].text.trim();\n  const nextText = items2[i + 1].text.trim();\n  if (prevText === \"Beginn des Texture-Designs\" && curText === \"\" && nextText === \"\") {\n    items2[i].insertText(\"6.06.2022: \", Word.InsertLocation.replace);\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1 -------------------------------------------------------------\n# The paragraph \"2.06.2022: Implementierung von weitern Funktionen der\n# GameEngine:\" is currently split across two runs (\"2.06.2022\" and\n# \": Implementierung ...\"). Find/Replace the whole sentence so Word\n# collapses it back down into a single run.\n$mergedText = \"2.06.2022: Implementierung von weitern Funktionen der GameEngine:\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $mergedText\n$find.Replacement.Text = $mergedText\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\n# --- Change 2 -------------------------------------------------------------\n# Of the two empty trailing paragraphs, the first one (immediately after\n# \"Beginn des Texture-Designs\") receives the new text \"6.06.2022: \"; the\n# very last (empty) paragraph of the document is left untouched.\nfor ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {\n    $prev = $d.Paragraphs.Item($i).Range.Text.TrimEnd()\n    $cur = $d.Paragraphs.Item($i + 1)\n    $curText = $cur.Range.Text.TrimEnd()\n    $hasNext = ($i + 1) -lt $d.Paragraphs.Count\n    if ($hasNext) {\n        $nextText = $d.Paragraphs.Item($i + 2).Range.Text.TrimEnd()\n    } else {\n        $nextText = $null\n    }\n    if ($prev -eq \"Beginn des Texture-Designs\" -and $curText -eq \"\" -and $nextText -eq \"\") {\n        $cur.Range.InsertBefore(\"6.06.2022: \")\n        break\n    }\n}\n"}
